# Update the "Förändrad" (Changed) date column (C) for data rows 2-44
# from serial date 45175 (2023-09-06) to serial date 45177 (2023-09-08).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

for ($row = 2; $row -le 44; $row++) {
    $cell = $ws.Cells.Item($row, 3)  # Column C
    if ($cell.Value2 -eq 45175) {
        $cell.Value2 = 45177
    }
}
